$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.866.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.27"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.62"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.588.29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.555"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.01"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.874.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.94"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.48%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.37"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.45"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.417.83"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.04%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.853"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.70"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.769.84"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.45%  "
